$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3475954063816255
$ws.Range("D2").Value = -0.3610913323653295
$ws.Range("E2").Value = 0.6410990123960496
$ws.Range("F2").Value = 0.1805350581402326
$ws.Range("G2").Value = 0.006229848919395677
$ws.Range("H2").Value = -0.2250753963015852
$ws.Range("I2").Value = 0.02697015588062352
$ws.Range("J2").Value = -0.4385348421393686
$ws.Range("K2").Value = 0.1687535070140281
$ws.Range("L2").Value = 0.0218563114252457
$ws.Range("M2").Value = 0.00005356821427285709
$ws.Range("N2").Value = -0.07790104760419042
$ws.Range("O2").Value = -0.003926031704126817
$ws.Range("P2").Value = -0.04108029232116928
$ws.Range("Q2").Value = -0.05146551786207145
$ws.Range("R2").Value = 0.0536185184740739
$ws.Range("S2").Value = -0.01882442729770919
$ws.Range("T2").Value = -0.02999848799395198
$ws.Range("U2").Value = 0.2021136404545618
$ws.Range("V2").Value = -0.02054715418861676
$ws.Range("W2").Value = 0.05227048108192433
$ws.Range("X2").Value = -0.03909817239268957
$ws.Range("Y2").Value = -0.0806942747770991
$ws.Range("Z2").Value = -0.04372529490117961
$ws.Range("AA2").Value = -0.07290163560654242
$ws.Range("AB2").Value = -0.03013663254653019
$ws.Range("AC2").Value = 0.08782057528230112
$ws.Range("AD2").Value = 0.1926756507026028
$ws.Range("AE2").Value = 0.1069668918675675
$ws.Range("AF2").Value = -0.115059628238513
$ws.Range("AG2").Value = -0.1306914667658671
$ws.Range("AH2").Value = -0.01157822231288925
$ws.Range("AI2").Value = -0.02099931599726399
$ws.Range("AJ2").Value = 0.5568619394477577
$ws.Range("AK2").Value = 0.4853875415501662
$ws.Range("B3").Value = -0.3475954063816255
$ws.Range("D3").Value = 0.3487166908667635
$ws.Range("E3").Value = 0.02553581414325657
$ws.Range("F3").Value = -0.2948527954111816
$ws.Range("G3").Value = 0.01640329761319045
$ws.Range("H3").Value = -0.238628346513386
$ws.Range("I3").Value = 0.4403663374653499
$ws.Range("J3").Value = 0.05312152448609794
$ws.Range("K3").Value = -0.04860201840807363
$ws.Range("L3").Value = -0.01835556142224569
$ws.Range("M3").Value = -0.05225080100320401
$ws.Range("N3").Value = 0.08445268981075925
$ws.Range("O3").Value = -0.0838253273013092
$ws.Range("P3").Value = 0.03505358021432086
$ws.Range("Q3").Value = 0.1181475765903064
$ws.Range("R3").Value = -0.01468229872919492
$ws.Range("S3").Value = -0.2976750147000588
$ws.Range("T3").Value = -0.03318973275893104
$ws.Range("U3").Value = 0.5012014928059713
$ws.Range("V3").Value = 0.07263398653594615
$ws.Range("W3").Value = -0.05950103800415202
$ws.Range("X3").Value = -0.03357191028764115
$ws.Range("Y3").Value = -0.06172229488917955
$ws.Range("Z3").Value = 0.04418676874707499
$ws.Range("AA3").Value = 0.02402822411289645
$ws.Range("AB3").Value = 0.01291944367777471
$ws.Range("AC3").Value = 0.03095292381169525
$ws.Range("AD3").Value = -0.2916204624818499
$ws.Range("AE3").Value = -0.00299252397009588
$ws.Range("AF3").Value = 0.0535249180996724
$ws.Range("AG3").Value = -0.01598233592934372
$ws.Range("AH3").Value = 0.01317931671726687
$ws.Range("AI3").Value = 0.02548234192936772
$ws.Range("AJ3").Value = -0.420935667742671
$ws.Range("AK3").Value = -0.3020755763023052
$ws.Range("B4").Value = -0.3610913323653295
$ws.Range("C4").Value = 0.3487166908667635
$ws.Range("E4").Value = -0.3478680474721899
$ws.Range("F4").Value = 0.0572563570254281
$ws.Range("G4").Value = -0.01014061656246625
$ws.Range("H4").Value = -0.2275808143232573
$ws.Range("I4").Value = 0.02059957039828159
$ws.Range("J4").Value = 0.6418239112956452
$ws.Range("K4").Value = 0.02253935415741663
$ws.Range("L4").Value = 0.2069156756627027
$ws.Range("M4").Value = 0.02646663386653546
$ws.Range("N4").Value = -0.01931047724190897
$ws.Range("O4").Value = 0.2677215508862035
$ws.Range("P4").Value = 0.03383600734402938
$ws.Range("Q4").Value = 0.03317744470977884
$ws.Range("R4").Value = -0.01128014112056448
$ws.Range("S4").Value = 0.007858495433981736
$ws.Range("T4").Value = -0.05011402445609783
$ws.Range("U4").Value = 0.3558919515678063
$ws.Range("V4").Value = 0.2001787527150109
$ws.Range("W4").Value = -0.05075924303697215
$ws.Range("X4").Value = 0.02572647090588362
$ws.Range("Y4").Value = 0.2659693518774075
$ws.Range("Z4").Value = 0.0109247476989908
$ws.Range("AA4").Value = 0.06978392713570854
$ws.Range("AB4").Value = -0.02540890163560654
$ws.Range("AC4").Value = 0.03163596654386618
$ws.Range("AD4").Value = 0.05590246360985444
$ws.Range("AE4").Value = -0.02251343405373622
$ws.Range("AF4").Value = 0.1241576806307225
$ws.Range("AG4").Value = 0.1059957519830079
$ws.Range("AH4").Value = 0.05139946959787839
$ws.Range("AI4").Value = 0.002612554450217801
$ws.Range("AJ4").Value = -0.2796407665630662
$ws.Range("AK4").Value = 0.07428413713654855
$ws.Range("B5").Value = 0.6410990123960496
$ws.Range("C5").Value = 0.02553581414325657
$ws.Range("D5").Value = -0.3478680474721899
$ws.Range("F5").Value = -0.2307834831339325
$ws.Range("G5").Value = 0.007962367849471397
$ws.Range("H5").Value = -0.1230694202776811
$ws.Range("I5").Value = 0.278566234264937
$ws.Range("J5").Value = -0.4991819327277309
$ws.Range("K5").Value = 0.2785856263425054
$ws.Range("L5").Value = 0.02072101088404354
$ws.Range("M5").Value = 0.07314192456769827
$ws.Range("N5").Value = -0.02206952827811311
$ws.Range("O5").Value = -0.001053700214800859
$ws.Range("P5").Value = -0.001401605606422426
$ws.Range("Q5").Value = 0.02152856611426446
$ws.Range("R5").Value = -0.05153444613778455
$ws.Range("S5").Value = -0.1746565706262825
$ws.Range("T5").Value = -0.07947852591410366
$ws.Range("U5").Value = 0.09085841143364573
$ws.Range("V5").Value = 0.01899089196356785
$ws.Range("W5").Value = -0.015687518750075
$ws.Range("X5").Value = 0.001142884571538286
$ws.Range("Y5").Value = -0.09146023784095136
$ws.Range("Z5").Value = -0.001891687566750267
$ws.Range("AA5").Value = -0.05023786495145981
$ws.Range("AB5").Value = -0.02631188124752499
$ws.Range("AC5").Value = 0.03766488665954664
$ws.Range("AD5").Value = -0.2302870651482606
$ws.Range("AE5").Value = 0.08727596510386042
$ws.Range("AF5").Value = -0.0554747818991276
$ws.Range("AG5").Value = -0.0556126384505538
$ws.Range("AH5").Value = -0.06841918567674271
$ws.Range("AI5").Value = -0.02803000012000048
$ws.Range("AJ5").Value = 0.05332446929787719
$ws.Range("AK5").Value = 0.5040794883179532
$ws.Range("B6").Value = 0.1805350581402326
$ws.Range("C6").Value = -0.2948527954111816
$ws.Range("D6").Value = 0.0572563570254281
$ws.Range("E6").Value = -0.2307834831339325
$ws.Range("G6").Value = -0.03901359605438422
$ws.Range("H6").Value = -0.03914770859083436
$ws.Range("I6").Value = -0.05511074844299377
$ws.Range("J6").Value = -0.05668361873447494
$ws.Range("K6").Value = 0.009321637286549146
$ws.Range("L6").Value = 0.007117180468721875
$ws.Range("M6").Value = -0.03450042600170401
$ws.Range("N6").Value = 0.01192411169644679
$ws.Range("O6").Value = -0.003325069300277201
$ws.Range("P6").Value = -0.002661994647978592
$ws.Range("Q6").Value = -0.0713376933507734
$ws.Range("R6").Value = 0.07848828195312781
$ws.Range("S6").Value = 0.003153612614450458
$ws.Range("T6").Value = -0.01942529370117481
$ws.Range("U6").Value = -0.06557853431413725
$ws.Range("V6").Value = 0.03911612446449786
$ws.Range("W6").Value = 0.09336392545570182
$ws.Range("X6").Value = 0.04664841059364237
$ws.Range("Y6").Value = 0.003374221496885987
$ws.Range("Z6").Value = 0.008808419233676935
$ws.Range("AA6").Value = -0.05389231156924628
$ws.Range("AB6").Value = -0.08107779631118524
$ws.Range("AC6").Value = 0.004654770619082476
$ws.Range("AD6").Value = 0.989935991743967
$ws.Range("AE6").Value = 0.008488641954567819
$ws.Range("AF6").Value = 0.02474141096564386
$ws.Range("AG6").Value = 0.01600441601766407
$ws.Range("AH6").Value = 0.01823796095184381
$ws.Range("AI6").Value = 0.0499750159000636
$ws.Range("AJ6").Value = -0.005805431221724887
$ws.Range("AK6").Value = 0.08772063888255553
$ws.Range("B7").Value = 0.006229848919395677
$ws.Range("C7").Value = 0.01640329761319045
$ws.Range("D7").Value = -0.01014061656246625
$ws.Range("E7").Value = 0.007962367849471397
$ws.Range("F7").Value = -0.03901359605438422
$ws.Range("B8").Value = -0.2250753963015852
$ws.Range("C8").Value = -0.238628346513386
$ws.Range("D8").Value = -0.2275808143232573
$ws.Range("E8").Value = -0.1230694202776811
$ws.Range("F8").Value = -0.03914770859083436
$ws.Range("B9").Value = 0.02697015588062352
$ws.Range("C9").Value = 0.4403663374653499
$ws.Range("D9").Value = 0.02059957039828159
$ws.Range("E9").Value = 0.278566234264937
$ws.Range("F9").Value = -0.05511074844299377
$ws.Range("B10").Value = -0.4385348421393686
$ws.Range("C10").Value = 0.05312152448609794
$ws.Range("D10").Value = 0.6418239112956452
$ws.Range("E10").Value = -0.4991819327277309
$ws.Range("F10").Value = -0.05668361873447494
$ws.Range("B11").Value = 0.1687535070140281
$ws.Range("C11").Value = -0.04860201840807363
$ws.Range("D11").Value = 0.02253935415741663
$ws.Range("E11").Value = 0.2785856263425054
$ws.Range("F11").Value = 0.009321637286549146
$ws.Range("B12").Value = 0.0218563114252457
$ws.Range("C12").Value = -0.01835556142224569
$ws.Range("D12").Value = 0.2069156756627027
$ws.Range("E12").Value = 0.02072101088404354
$ws.Range("F12").Value = 0.007117180468721875
$ws.Range("B13").Value = 0.00005356821427285709
$ws.Range("C13").Value = -0.05225080100320401
$ws.Range("D13").Value = 0.02646663386653546
$ws.Range("E13").Value = 0.07314192456769827
$ws.Range("F13").Value = -0.03450042600170401
$ws.Range("B14").Value = -0.07790104760419042
$ws.Range("C14").Value = 0.08445268981075925
$ws.Range("D14").Value = -0.01931047724190897
$ws.Range("E14").Value = -0.02206952827811311
$ws.Range("F14").Value = 0.01192411169644679
$ws.Range("B15").Value = -0.003926031704126817
$ws.Range("C15").Value = -0.0838253273013092
$ws.Range("D15").Value = 0.2677215508862035
$ws.Range("E15").Value = -0.001053700214800859
$ws.Range("F15").Value = -0.003325069300277201
$ws.Range("B16").Value = -0.04108029232116928
$ws.Range("C16").Value = 0.03505358021432086
$ws.Range("D16").Value = 0.03383600734402938
$ws.Range("E16").Value = -0.001401605606422426
$ws.Range("F16").Value = -0.002661994647978592
$ws.Range("B17").Value = -0.05146551786207145
$ws.Range("C17").Value = 0.1181475765903064
$ws.Range("D17").Value = 0.03317744470977884
$ws.Range("E17").Value = 0.02152856611426446
$ws.Range("F17").Value = -0.0713376933507734
$ws.Range("B18").Value = 0.0536185184740739
$ws.Range("C18").Value = -0.01468229872919492
$ws.Range("D18").Value = -0.01128014112056448
$ws.Range("E18").Value = -0.05153444613778455
$ws.Range("F18").Value = 0.07848828195312781
$ws.Range("B19").Value = -0.01882442729770919
$ws.Range("C19").Value = -0.2976750147000588
$ws.Range("D19").Value = 0.007858495433981736
$ws.Range("E19").Value = -0.1746565706262825
$ws.Range("F19").Value = 0.003153612614450458
$ws.Range("B20").Value = -0.02999848799395198
$ws.Range("C20").Value = -0.03318973275893104
$ws.Range("D20").Value = -0.05011402445609783
$ws.Range("E20").Value = -0.07947852591410366
$ws.Range("F20").Value = -0.01942529370117481
$ws.Range("B21").Value = 0.2021136404545618
$ws.Range("C21").Value = 0.5012014928059713
$ws.Range("D21").Value = 0.3558919515678063
$ws.Range("E21").Value = 0.09085841143364573
$ws.Range("F21").Value = -0.06557853431413725
$ws.Range("B22").Value = -0.02054715418861676
$ws.Range("C22").Value = 0.07263398653594615
$ws.Range("D22").Value = 0.2001787527150109
$ws.Range("E22").Value = 0.01899089196356785
$ws.Range("F22").Value = 0.03911612446449786
$ws.Range("B23").Value = 0.05227048108192433
$ws.Range("C23").Value = -0.05950103800415202
$ws.Range("D23").Value = -0.05075924303697215
$ws.Range("E23").Value = -0.015687518750075
$ws.Range("F23").Value = 0.09336392545570182
$ws.Range("B24").Value = -0.03909817239268957
$ws.Range("C24").Value = -0.03357191028764115
$ws.Range("D24").Value = 0.02572647090588362
$ws.Range("E24").Value = 0.001142884571538286
$ws.Range("F24").Value = 0.04664841059364237
$ws.Range("B25").Value = -0.0806942747770991
$ws.Range("C25").Value = -0.06172229488917955
$ws.Range("D25").Value = 0.2659693518774075
$ws.Range("E25").Value = -0.09146023784095136
$ws.Range("F25").Value = 0.003374221496885987
$ws.Range("B26").Value = -0.04372529490117961
$ws.Range("C26").Value = 0.04418676874707499
$ws.Range("D26").Value = 0.0109247476989908
$ws.Range("E26").Value = -0.001891687566750267
$ws.Range("F26").Value = 0.008808419233676935
$ws.Range("B27").Value = -0.07290163560654242
$ws.Range("C27").Value = 0.02402822411289645
$ws.Range("D27").Value = 0.06978392713570854
$ws.Range("E27").Value = -0.05023786495145981
$ws.Range("F27").Value = -0.05389231156924628
$ws.Range("B28").Value = -0.03013663254653019
$ws.Range("C28").Value = 0.01291944367777471
$ws.Range("D28").Value = -0.02540890163560654
$ws.Range("E28").Value = -0.02631188124752499
$ws.Range("F28").Value = -0.08107779631118524
$ws.Range("B29").Value = 0.08782057528230112
$ws.Range("C29").Value = 0.03095292381169525
$ws.Range("D29").Value = 0.03163596654386618
$ws.Range("E29").Value = 0.03766488665954664
$ws.Range("F29").Value = 0.004654770619082476
$ws.Range("B30").Value = 0.1926756507026028
$ws.Range("C30").Value = -0.2916204624818499
$ws.Range("D30").Value = 0.05590246360985444
$ws.Range("E30").Value = -0.2302870651482606
$ws.Range("F30").Value = 0.989935991743967
$ws.Range("B31").Value = 0.1069668918675675
$ws.Range("C31").Value = -0.00299252397009588
$ws.Range("D31").Value = -0.02251343405373622
$ws.Range("E31").Value = 0.08727596510386042
$ws.Range("F31").Value = 0.008488641954567819
$ws.Range("B32").Value = -0.115059628238513
$ws.Range("C32").Value = 0.0535249180996724
$ws.Range("D32").Value = 0.1241576806307225
$ws.Range("E32").Value = -0.0554747818991276
$ws.Range("F32").Value = 0.02474141096564386
$ws.Range("B33").Value = -0.1306914667658671
$ws.Range("C33").Value = -0.01598233592934372
$ws.Range("D33").Value = 0.1059957519830079
$ws.Range("E33").Value = -0.0556126384505538
$ws.Range("F33").Value = 0.01600441601766407
$ws.Range("B34").Value = -0.01157822231288925
$ws.Range("C34").Value = 0.01317931671726687
$ws.Range("D34").Value = 0.05139946959787839
$ws.Range("E34").Value = -0.06841918567674271
$ws.Range("F34").Value = 0.01823796095184381
$ws.Range("B35").Value = -0.02099931599726399
$ws.Range("C35").Value = 0.02548234192936772
$ws.Range("D35").Value = 0.002612554450217801
$ws.Range("E35").Value = -0.02803000012000048
$ws.Range("F35").Value = 0.0499750159000636
$ws.Range("B36").Value = 0.5568619394477577
$ws.Range("C36").Value = -0.420935667742671
$ws.Range("D36").Value = -0.2796407665630662
$ws.Range("E36").Value = 0.05332446929787719
$ws.Range("F36").Value = -0.005805431221724887
$ws.Range("B37").Value = 0.4853875415501662
$ws.Range("C37").Value = -0.3020755763023052
$ws.Range("D37").Value = 0.07428413713654855
$ws.Range("E37").Value = 0.5040794883179532
$ws.Range("F37").Value = 0.08772063888255553
